$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.451.66"
$ws.Range("E2").Value = "  +1.69%  "
$ws.Range("D3").Value = "2.158.87"
$ws.Range("E3").Value = "  +2.59%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.82"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.98"
$ws.Range("E7").Value = "  +2.44%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  +2.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0855"
$ws.Range("E10").Value = "  +1.14%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.99"
$ws.Range("E12").Value = "  +0.99%  "
$ws.Range("D13").Value = "2.478.76"
$ws.Range("E13").Value = "  +2.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.10"
$ws.Range("E14").Value = "  +0.12%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.813"
$ws.Range("E15").Value = "  +0.65%  "
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").Value = "2.157.73"
$ws.Range("E17").Value = "  +2.70%  "
$ws.Range("D18").Value = "39.381.80"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.86"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("E20").Value = "  +0.66%  "
$ws.Range("D21").Value = "0.0₃0850"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.43"
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.35"
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("E25").Value = "  -3.89%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "172.24"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.50"
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.143"
$ws.Range("E28").Value = "  +3.59%  "
$ws.Range("E29").Value = "  +2.52%  "
$ws.Range("E30").Value = "  +0.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.63"
$ws.Range("E31").Value = "  +5.01%  "
$ws.Range("E32").Value = "  +0.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.60"
$ws.Range("E33").Value = "  +1.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "7.14"
$ws.Range("E34").Value = "  +8.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.73"
$ws.Range("E35").Value = "  -0.75%  "
$ws.Range("E36").Value = "  -1.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.42"
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("E38").Value = "  -0.27%  "
$ws.Range("E39").Value = "  +0.12%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "103.61"
$ws.Range("E40").Value = "  +1.92%  "
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.65"
$ws.Range("E42").Value = "  -3.76%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.530.99"
$ws.Range("E43").Value = "  -0.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.19"
$ws.Range("E44").Value = "  +4.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.34"
$ws.Range("E45").Value = "  +5.43%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0930"
$ws.Range("E46").Value = "  +1.87%  "
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.10"
$ws.Range("E48").Value = "  +4.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.69"
$ws.Range("E49").Value = "  -1.29%  "
$ws.Range("D50").Value = "2.362.22"
$ws.Range("E50").Value = "  +2.85%  "
$ws.Range("E51").Value = "  -0.16%  "
